# Apply "Add files via upload" changes to database.xlsx
# Sheets: "user" (sheet1), "driver" (sheet2), "admin" (sheet3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "user" sheet: update a few existing cells and append new rows 11-23
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("user")

# Row 7: id/name text stays logically the same ("new" / "krish")
$wsUser.Range("B7").Value = "new"
$wsUser.Range("D7").Value = "krish"

# Row 8
$wsUser.Range("B8").Value = "1JB0"

# Row 9: id changes, dept becomes MECH, usn cell (I9) is cleared
$wsUser.Range("B9").Value = "1JB1"
$wsUser.Range("G9").Value = "MECH"
$wsUser.Range("I9").ClearContents()

# Row 10: id changes, dept becomes ISE
$wsUser.Range("B10").Value = "1JB2"
$wsUser.Range("G10").Value = "ISE"
$wsUser.Range("I10").Value = "usn"

# New rows 11-14 (full data: id, password, name, dept, sem, usn)
$wsUser.Range("A11").Value = 10
$wsUser.Range("B11").Value = "1JB10"
$wsUser.Range("C11").Value = 9052004
$wsUser.Range("D11").Value = "Jayesh"
$wsUser.Range("G11").Value = "CSE"
$wsUser.Range("H11").Value = 5
$wsUser.Range("I11").Value = "1JB23CS070"

$wsUser.Range("A12").Value = 11
$wsUser.Range("B12").Value = "1JB11"
$wsUser.Range("C12").Value = 1012000
$wsUser.Range("D12").Value = "Saurya"
$wsUser.Range("G12").Value = "CSE"
$wsUser.Range("H12").Value = 5
$wsUser.Range("I12").Value = "1JB23CS080"

$wsUser.Range("A13").Value = 12
$wsUser.Range("B13").Value = "1JB12"
$wsUser.Range("C13").Value = 1012000
$wsUser.Range("D13").Value = "Krish"
$wsUser.Range("G13").Value = "CSE"
$wsUser.Range("H13").Value = 5
$wsUser.Range("I13").Value = "1JB23CS073"

$wsUser.Range("A14").Value = 13
$wsUser.Range("B14").Value = "1JB13"
$wsUser.Range("C14").Value = 1012000
$wsUser.Range("D14").Value = "John S Mark"
$wsUser.Range("G14").Value = "CSE"
$wsUser.Range("H14").Value = 5
$wsUser.Range("I14").Value = "1JB23CS075"

# New rows 15-23: just id + password (A,B,C only)
$wsUser.Range("A15").Value = 14
$wsUser.Range("B15").Value = "1JB14"
$wsUser.Range("C15").Value = 1012000

$wsUser.Range("A16").Value = 15
$wsUser.Range("B16").Value = "1JB15"
$wsUser.Range("C16").Value = 1012000

$wsUser.Range("A17").Value = 16
$wsUser.Range("B17").Value = "1JB16"
$wsUser.Range("C17").Value = 1012000

$wsUser.Range("A18").Value = 17
$wsUser.Range("B18").Value = "1JB17"
$wsUser.Range("C18").Value = 1012000

$wsUser.Range("A19").Value = 18
$wsUser.Range("B19").Value = "1JB18"
$wsUser.Range("C19").Value = 1012000

$wsUser.Range("A20").Value = 19
$wsUser.Range("B20").Value = "1JB19"
$wsUser.Range("C20").Value = 1012000

$wsUser.Range("A21").Value = 20
$wsUser.Range("B21").Value = "1JB20"
$wsUser.Range("C21").Value = 1012000

$wsUser.Range("A22").Value = 21
$wsUser.Range("B22").Value = "1JB21"
$wsUser.Range("C22").Value = 1012000

$wsUser.Range("A23").Value = 22
$wsUser.Range("B23").Value = "1JB22"
$wsUser.Range("C23").Value = 1012000

$wsUser.Activate()
$wsUser.Range("I15").Select()

# ---------------------------------------------------------------------------
# 2) "driver" sheet: the driver named "driverr" (row 2) is removed and every
#    later driver shifts up one row, leaving row 6 blank (style-only on B6)
# ---------------------------------------------------------------------------
$wsDriver = $wb.Worksheets.Item("driver")
$wsDriver.Range("B2").Value = "umesh"
$wsDriver.Range("B3").Value = "pradeep"
$wsDriver.Range("B4").Value = "paramesh"
$wsDriver.Range("B5").Value = "chandru"
$wsDriver.Range("A6:C6").ClearContents()

$wsDriver.Activate()
$wsDriver.Range("C9").Select()

# ---------------------------------------------------------------------------
# 3) "admin" sheet: rename the admin account name from "adminn" to "admin"
# ---------------------------------------------------------------------------
$wsAdmin = $wb.Worksheets.Item("admin")
$wsAdmin.Range("B2").Value = "admin"

$wsAdmin.Activate()
$wsAdmin.Range("B2").Select()

# Re-activate the "user" sheet, which is the tab selected in the saved file
$wsUser.Activate()

"done"
